# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to currentAveragePrice / LevePrice / LeveProfit columns
# across multiple sheets, per the authoritative diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 8512
$ws.Range("I34").Value = 8512
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 8512
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -8309
$ws.Range("N34").ClearContents()
# Row 36
$ws.Range("H36").Value = 8512
$ws.Range("I36").Value = 8512
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 8512
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -7797
$ws.Range("N36").ClearContents()
# Row 41
$ws.Range("H41").Value = 676.1111
$ws.Range("J41").Value = 737.8
$ws.Range("L41").Value = 737.8
$ws.Range("N41").Value = -1617.8
# Row 80
$ws.Range("H80").Value = 2874.2
$ws.Range("I80").Value = 1374.25
$ws.Range("K80").Value = 4122.75
$ws.Range("M80").Value = -3124.75
# Row 83
$ws.Range("H83").Value = 2874.2
$ws.Range("I83").Value = 1374.25
$ws.Range("K83").Value = 12368.25
$ws.Range("M83").Value = -7376.25
# Row 92
$ws.Range("H92").Value = 1422.2778
$ws.Range("I92").Value = 1450.5
$ws.Range("J92").Value = 1365.8334
$ws.Range("K92").Value = 1450.5
$ws.Range("L92").Value = 1365.8334
$ws.Range("M92").Value = -202.5
$ws.Range("N92").Value = -3861.8334
# Row 96
$ws.Range("H96").Value = 2026.1666
$ws.Range("I96").Value = 431.6
$ws.Range("J96").Value = 9999
$ws.Range("K96").Value = 1294.8
$ws.Range("L96").Value = 29997
$ws.Range("M96").Value = 78.19999999999982
$ws.Range("N96").Value = -32743
# Row 98
$ws.Range("H98").Value = 2614
$ws.Range("I98").Value = 1008
$ws.Range("K98").Value = 1008
$ws.Range("M98").Value = 490
# Row 116
$ws.Range("H116").Value = 9156
$ws.Range("I116").Value = 9926.666999999999
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 9926.666999999999
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = -6484.666999999999
$ws.Range("N116").Value = -14884
# Row 118
$ws.Range("H118").Value = 2451.8684
$ws.Range("I118").Value = 396.5
$ws.Range("J118").Value = 2999.9666
$ws.Range("K118").Value = 1189.5
$ws.Range("L118").Value = 8999.899800000001
$ws.Range("M118").Value = 467.5
$ws.Range("N118").Value = -12313.8998
# Row 122
$ws.Range("H122").Value = 2614
$ws.Range("I122").Value = 1008
$ws.Range("K122").Value = 3024
$ws.Range("M122").Value = -574
# Row 138
$ws.Range("H138").Value = 5546.1724
$ws.Range("J138").Value = 6926.476
$ws.Range("L138").Value = 20779.428
$ws.Range("N138").Value = -31059.428

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1698.75
$ws.Range("I2").Value = 1798.3334
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1798.3334
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -1685.3334
$ws.Range("N2").Value = -1626
# Row 32
$ws.Range("H32").Value = 12223.679
$ws.Range("I32").Value = 10690.52
$ws.Range("K32").Value = 10690.52
$ws.Range("M32").Value = -10403.52
# Row 74
$ws.Range("H74").Value = 11331.167
$ws.Range("I74").Value = 9993.429
$ws.Range("J74").Value = 13204
$ws.Range("K74").Value = 9993.429
$ws.Range("L74").Value = 13204
$ws.Range("M74").Value = -9119.429
$ws.Range("N74").Value = -14952
# Row 77
$ws.Range("H77").Value = 11331.167
$ws.Range("I77").Value = 9993.429
$ws.Range("J77").Value = 13204
$ws.Range("K77").Value = 49967.145
$ws.Range("L77").Value = 66020
$ws.Range("M77").Value = -45599.145
$ws.Range("N77").Value = -74756
# Row 116
$ws.Range("H116").Value = 1698.75
$ws.Range("I116").Value = 1798.3334
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 1798.3334
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = 495.6666
$ws.Range("N116").Value = -5988
# Row 132
$ws.Range("H132").Value = 3241.25
$ws.Range("I132").Value = 2669.5454
$ws.Range("K132").Value = 8008.6362
$ws.Range("M132").Value = -5478.6362

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1698.75
$ws.Range("I3").Value = 1798.3334
$ws.Range("J3").Value = 1400
$ws.Range("K3").Value = 1798.3334
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = -1684.3334
$ws.Range("N3").Value = -1628
# Row 99
$ws.Range("H99").Value = 21241.8
$ws.Range("I99").Value = 26052.25
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 26052.25
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -24554.25
$ws.Range("N99").Value = -4996
# Row 107
$ws.Range("H107").Value = 1848.75
$ws.Range("I107").Value = 1766.6666
$ws.Range("K107").Value = 1766.6666
$ws.Range("M107").Value = 153.3334

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2666
$ws.Range("I16").Value = 2499.5
$ws.Range("J16").Value = 2999
$ws.Range("K16").Value = 2499.5
$ws.Range("L16").Value = 2999
$ws.Range("M16").Value = -2212.5
$ws.Range("N16").Value = -3573
# Row 31
$ws.Range("H31").Value = 2266
$ws.Range("I31").Value = 2036.2222
$ws.Range("K31").Value = 2036.2222
$ws.Range("M31").Value = -1741.2222
# Row 34
$ws.Range("H34").Value = 2266
$ws.Range("I34").Value = 2036.2222
$ws.Range("K34").Value = 2036.2222
$ws.Range("M34").Value = -1834.2222
# Row 113
$ws.Range("H113").Value = 2666
$ws.Range("I113").Value = 2499.5
$ws.Range("J113").Value = 2999
$ws.Range("K113").Value = 2499.5
$ws.Range("L113").Value = 2999
$ws.Range("M113").Value = -329.5
$ws.Range("N113").Value = -7339
# Row 132
$ws.Range("H132").Value = 2912.5454
$ws.Range("I132").Value = 1994.1666
$ws.Range("K132").Value = 5982.4998
$ws.Range("M132").Value = -3452.4998

$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 1668000
$ws.Range("J32").Value = 1668000
$ws.Range("L32").Value = 5004000
$ws.Range("N32").Value = -5004566
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# Row 86
$ws.Range("H86").Value = 400
$ws.Range("I86").Value = 400
$ws.Range("K86").Value = 1200
$ws.Range("M86").Value = -14
# Row 89
$ws.Range("H89").Value = 400
$ws.Range("I89").Value = 400
$ws.Range("K89").Value = 3600
$ws.Range("M89").Value = 2328
# Row 107
$ws.Range("H107").Value = 648.2632
$ws.Range("J107").Value = 671.58826
$ws.Range("L107").Value = 2014.76478
$ws.Range("N107").Value = -5854.76478
# Row 113
$ws.Range("H113").Value = 1943.2727
$ws.Range("J113").Value = 1998.375
$ws.Range("L113").Value = 5995.125
$ws.Range("N113").Value = -10335.125
# Row 128
$ws.Range("H128").Value = 790000
$ws.Range("I128").Value = 790000
$ws.Range("K128").Value = 2370000
$ws.Range("M128").Value = -2365020

$ws = $wb.Worksheets.Item("GSM")
# Row 74
$ws.Range("H74").Value = 49131
$ws.Range("J74").Value = 49131
$ws.Range("L74").Value = 49131
$ws.Range("N74").Value = -51003
# Row 77
$ws.Range("H77").Value = 49131
$ws.Range("J77").Value = 49131
$ws.Range("L77").Value = 147393
$ws.Range("N77").Value = -156753
# Row 80
$ws.Range("H80").Value = 2972.75
$ws.Range("I80").Value = 1997
$ws.Range("J80").Value = 3948.5
$ws.Range("K80").Value = 1997
$ws.Range("L80").Value = 3948.5
$ws.Range("M80").Value = -999
$ws.Range("N80").Value = -5944.5
# Row 83
$ws.Range("H83").Value = 2972.75
$ws.Range("I83").Value = 1997
$ws.Range("J83").Value = 3948.5
$ws.Range("K83").Value = 9985
$ws.Range("L83").Value = 19742.5
$ws.Range("M83").Value = -4993
$ws.Range("N83").Value = -29726.5
# Row 122
$ws.Range("H122").Value = 3858.762
$ws.Range("I122").Value = 2768.6667
$ws.Range("J122").Value = 10399.333
$ws.Range("K122").Value = 8306.000100000001
$ws.Range("L122").Value = 31197.999
$ws.Range("M122").Value = -5856.000100000001
$ws.Range("N122").Value = -36097.999
# Row 132
$ws.Range("H132").Value = 5437.5557
$ws.Range("I132").Value = 4134.143
$ws.Range("J132").Value = 9999.5
$ws.Range("K132").Value = 12402.429
$ws.Range("L132").Value = 29998.5
$ws.Range("M132").Value = -9872.429
$ws.Range("N132").Value = -35058.5

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1791.5
$ws.Range("I61").Value = 1772.4
$ws.Range("J61").Value = 1887
$ws.Range("K61").Value = 1772.4
$ws.Range("L61").Value = 1887
$ws.Range("M61").Value = -1570.4
$ws.Range("N61").Value = -2291
# Row 113
$ws.Range("H113").Value = 1791.5
$ws.Range("I113").Value = 1772.4
$ws.Range("J113").Value = 1887
$ws.Range("K113").Value = 1772.4
$ws.Range("L113").Value = 1887
$ws.Range("M113").Value = 397.5999999999999
$ws.Range("N113").Value = -6227

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1442.48
$ws.Range("I126").Value = 780.4286
$ws.Range("K126").Value = 2341.2858
$ws.Range("M126").Value = 128.7142000000003
# Row 136
$ws.Range("H136").Value = 1033.875
$ws.Range("I136").Value = 1033.875
$ws.Range("M136").Value = -551.625

